$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("ById")
$ws3 = $wb.Worksheets.Item("ByGeoCoords")
$ws4 = $wb.Worksheets.Item("ByZipCode")

# ---------------------------------------------------------------------
# Step 1: seed brand-new shared strings in the exact order they need to
# land in sharedStrings.xml (26 Lat .. 35 ZZ).
# ---------------------------------------------------------------------
$ws3.Range("D1").Value = "Lat"
$ws3.Range("E1").Value = "Lon"
$ws3.Range("B2").Value = "Valid Geo Coords"
$ws3.Range("B3").Value = "Invalid Geo Coords"
$ws4.Range("E1").Value = "Country Code"
$ws4.Range("D1").Value = "Zip Code"
$ws4.Range("E2").Value = "us"
$ws4.Range("M2").Value = "Mountain View"
$ws4.Range("D3").Value = "XXXXX"
$ws4.Range("E3").Value = "ZZ"

# ---------------------------------------------------------------------
# Step 2: finish populating ByGeoCoords (sheet3)
# ---------------------------------------------------------------------
$ws3.Range("A1").Value = "TestCaseId"
$ws3.Range("B1").Value = "Description"
$ws3.Range("C1").Value = "Type"
$ws3.Range("F1").Value = "Format"
$ws3.Range("G1").Value = "Search Accuracy"
$ws3.Range("H1").Value = "Unit Format"
$ws3.Range("I1").Value = "Language"
$ws3.Range("J1").Value = "Expected Code"
$ws3.Range("K1").Value = "Expected Message"
$ws3.Range("L1").Value = "Expected Id"
$ws3.Range("M1").Value = "Expected Name"

$ws3.Range("A2").Value = 3001
$ws3.Range("C2").Value = "CORRECT"
$ws3.Range("D2").Value = -16.92
$ws3.Range("E2").Value = 145.77
$ws3.Range("F2").Value = "JSON"
$ws3.Range("G2").Value = "DEFAULT"
$ws3.Range("H2").Value = "STANDARD"
$ws3.Range("I2").Value = "DEFAULT"
$ws3.Range("J2").Value = 200
$ws3.Range("L2").Value = 2172797

$ws3.Range("A3").Formula = "=A2+1"
$ws3.Range("C3").Value = "INVALID_DATA"
$ws3.Range("D3").Value = 360
$ws3.Range("E3").Value = -1000
$ws3.Range("F3").Value = "JSON"
$ws3.Range("G3").Value = "DEFAULT"
$ws3.Range("H3").Value = "STANDARD"
$ws3.Range("I3").Value = "ENGLISH"
$ws3.Range("J3").Value = 404
$ws3.Range("K3").Value = "Error: Not found city"

# Re-use the existing named style (s="1", Arial Unicode MS / vertical-center)
# already present on sheet2 (ById) so no duplicate style/font gets created.
$ws2.Range("L2").Copy()
$ws3.Range("M2").PasteSpecial(-4122)
$ws3.Range("M2").Value = "Cairns"
$ws3.Range("M3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 3: finish populating ByZipCode (sheet4)
# ---------------------------------------------------------------------
$ws4.Range("A1").Value = "TestCaseId"
$ws4.Range("B1").Value = "Description"
$ws4.Range("C1").Value = "Type"
$ws4.Range("F1").Value = "Format"
$ws4.Range("G1").Value = "Search Accuracy"
$ws4.Range("H1").Value = "Unit Format"
$ws4.Range("I1").Value = "Language"
$ws4.Range("J1").Value = "Expected Code"
$ws4.Range("K1").Value = "Expected Message"
$ws4.Range("L1").Value = "Expected Id"
$ws4.Range("M1").Value = "Expected Name"

$ws4.Range("A2").Value = 3001
$ws4.Range("B2").Value = "Valid Geo Coords"
$ws4.Range("C2").Value = "CORRECT"
$ws4.Range("D2").Value = 94040
$ws4.Range("F2").Value = "JSON"
$ws4.Range("G2").Value = "DEFAULT"
$ws4.Range("H2").Value = "STANDARD"
$ws4.Range("I2").Value = "DEFAULT"
$ws4.Range("J2").Value = 200
$ws4.Range("L2").Value = 5375480

$ws4.Range("A3").Formula = "=A2+1"
$ws4.Range("B3").Value = "Invalid Geo Coords"
$ws4.Range("C3").Value = "INVALID_DATA"
$ws4.Range("F3").Value = "JSON"
$ws4.Range("G3").Value = "DEFAULT"
$ws4.Range("H3").Value = "STANDARD"
$ws4.Range("I3").Value = "ENGLISH"
$ws4.Range("J3").Value = 404
$ws4.Range("K3").Value = "Error: Not found city"

$ws2.Range("L2").Copy()
$ws4.Range("M2").PasteSpecial(-4122)
$ws4.Range("M2").Value = "Mountain View"
$ws4.Range("M3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 4: page setup (matches sheet1/sheet2's portrait orientation)
# ---------------------------------------------------------------------
$ws3.PageSetup.Orientation = 1
$ws4.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Step 5: selections / active sheet / active tab
# ---------------------------------------------------------------------
$ws2.Range("A1:L3").Select() | Out-Null
$ws3.Range("A1:M3").Select() | Out-Null
$ws4.Select()
$ws4.Range("K10").Select() | Out-Null
